# Daily attendance processing - rotate the "Recorded By" (column G) list
# for every data row: move the first name/email in the comma-separated
# list to the end (left-rotate by one). Cells with only a single
# recorder are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Text

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }
    if ($val -notmatch ",") {
        continue
    }

    $parts = $val -split ", "
    if ($parts.Count -lt 2) {
        continue
    }

    $first = $parts[0]
    $rest = $parts[1..($parts.Count - 1)]
    $newParts = $rest + @($first)
    $newVal = $newParts -join ", "

    $cell.Value = $newVal
}
